$d = $word.ActiveDocument

# 1. "Also from the t-value table " -> Null hypothesis statement
$d.Content.Find.Execute(
    "Also from the t-value table ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The Null hypothesis is taken as that there is no correlation between the series i.e. rho_population = 0",
    2) | Out-Null

# 2. Insert a new blank paragraph right before the "The observed t-value" paragraph
#    (splits the paragraph that used to directly precede it into two blank paragraphs).
$pTarget = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("The observed t-value ")) {
        $pTarget = $p
        break
    }
}
$full = $pTarget.Range
$insertPos = $full.End - 1
$markRange = $d.Range($insertPos, $insertPos)
$markRange.InsertBefore([char]13)

# 3. Rewrite the "The observed t-value ..." paragraph:
#    - lowercase "the observed t-value "
#    - split the trailing sentence, inserting ", from " + a HYPERLINK field
#      ("t-value table" pointing at wikipedia) before the remaining sentence.
$pTarget2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("The observed t-value ")) {
        $pTarget2 = $p
        break
    }
}
$full2 = $pTarget2.Range
$textRange = $d.Range($full2.Start, $full2.End - 1)
$textRange.Delete()
$insertPos2 = $textRange.Start
$ins = $d.Range($insertPos2, $insertPos2)
$xml = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">the observed t-value </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/></w:rPr><w:t>233.344275335</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, from </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink.0"/><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin" w:fldLock="0"/></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink.0"/><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "https://en.wikipedia.org/wiki/Student%27s_t-distribution#Table_of_selected_values"</w:instrText></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink.0"/><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate" w:fldLock="0"/></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink.0"/><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>t-value table</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end" w:fldLock="0"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> the t-value for 99.9% significance for a two sided test should be &gt;= 3.373 for 120 degrees of freedom, however we have 4999998 degrees of freedom, since the required t-value for significance decreases with increase in degrees of freedom, we can surely say that our correlation is significant.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$ins.InsertXML($xml)

# 4. "Thus we can reject the null hypothesis stating that rho is equal to 0." -> shortened
$d.Content.Find.Execute(
    "Thus we can reject the null hypothesis stating that rho is equal to 0.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Thus we can reject the null hypothesis",
    2) | Out-Null

# 5. Append "(-0.1055)" to the "not strong" sentence
$d.Content.Find.Execute(
    "However the relationship between the two series is not strong as the value is quite close to 0", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "However the relationship between the two series is not strong as the value is quite close to 0 (-0.1055)",
    2) | Out-Null

# 6. Remove the leading "However " from the "given more time" sentence
$d.Content.Find.Execute(
    "However given more time, would like to see in the factors such as cross correlation, or figuring out a non linear relationship between the series. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Given more time, would like to see in the factors such as cross correlation, or figuring out a non linear relationship between the series. ",
    2) | Out-Null

# 7. Footer page-number field cached value: 2 -> 1
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers(1)
    $ftr.Range.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2) | Out-Null
}

Write-Output "edits applied"
